# Update Name of Algo
# Applies updated imputed KNN result values to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value  = -7.944
$ws.Range("D6").Value  = -7.858
$ws.Range("C7").Value  = -12.659
$ws.Range("D7").Value  = -8.026
$ws.Range("A8").Value  = -21.181
$ws.Range("D8").Value  = -7.607000000000001
$ws.Range("D9").Value  = -8.032999999999999
$ws.Range("A10").Value = -20.945
$ws.Range("D10").Value = -7.696000000000001
$ws.Range("A12").Value = -21.649
$ws.Range("D12").Value = -8.364999999999998
$ws.Range("B13").Value = 6.476999999999999
$ws.Range("A18").Value = -21.649
$ws.Range("C20").Value = -13.041
$ws.Range("A25").Value = -21.534
